$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append (row 195..198).
$data = @(
    @("08-10-2021", 1.78, 2.7, 3.55, 4.24, -0.18),
    @("12-10-2021", 1.89, 2.9, 3.73, 4.49, -0.06),
    @("13-10-2021", 2.15, 2.97, 3.8, 4.67, 0.01),
    @("14-10-2021", 2.65, 3.22, 4.01, 4.92, 0.05)
)

$startRow = 195
$helper = $ws.Range("A300")

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]

    # Column A holds a date-formatted label that must stay plain text (it
    # matches existing "dd-mm-yyyy" shared strings elsewhere in the sheet).
    # Some of these labels (day <= 12) are ambiguous and would otherwise be
    # auto-converted to a date serial by a direct .Value assignment, so we
    # stage the text in a scratch cell forced to Text format, then copy only
    # the resulting value (not the format) into place.
    $helper.NumberFormat = "@"
    $helper.Value = $values[0]
    $helper.Copy()
    $ws.Range("A" + $row).PasteSpecial(-4163) | Out-Null

    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
}

# Remove the scratch helper cell/row entirely so it leaves no trace in the
# saved sheet (no stray cell, no dimension growth).
$ws.Rows("300:300").Delete()
